# Update "想去人数" (column F) values for specific events on the
# "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8498
$ws1.Range("F13").Value = 3677
$ws1.Range("F19").Value = 494
$ws1.Range("F22").Value = 1320
$ws1.Range("F24").Value = 434
$ws1.Range("F29").Value = 56
$ws1.Range("F33").Value = 654

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8498
$ws4.Range("F15").Value = 3677
$ws4.Range("F23").Value = 494
$ws4.Range("F27").Value = 1320
$ws4.Range("F29").Value = 434
$ws4.Range("F35").Value = 56
$ws4.Range("F39").Value = 654
